$d = $word.ActiveDocument

$pairs = @(
    @("347×8=2776", "994×2=1988"),
    @("294×9=2646", "297×7=2079"),
    @("348×9=3132", "311×8=2488"),
    @("418×4=1672", "864×6=5184"),
    @("267×9=2403", "391×8=3128"),
    @("601×9=5409", "806×3=2418"),
    @("158×9=1422", "296×6=1776"),
    @("178×2=356",  "177×3=531"),
    @("529×3=1587", "673×3=2019"),
    @("987×3=2961", "428×5=2140"),
    @("316×4=1264", "699×7=4893"),
    @("398×7=2786", "407×4=1628"),
    @("442×3=1326", "280×8=2240"),
    @("879×5=4395", "635×4=2540"),
    @("376×4=1504", "876×9=7884"),
    @("386×3=1158", "511×3=1533"),
    @("790×5=3950", "123×9=1107"),
    @("625×4=2500", "961×6=5766"),
    @("542×6=3252", "865×6=5190"),
    @("704×4=2816", "978×8=7824"),
    @("251×7=1757", "737×6=4422"),
    @("562×8=4496", "355×7=2485"),
    @("367×8=2936", "210×3=630"),
    @("716×3=2148", "739×8=5912"),
    @("445×6=2670", "418×9=3762")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
